# Update countries & provincias Spain
# Applies the 31-Mar-2020 11:20 data refresh to the "Pais" sheet:
#  - updates the "last updated" timestamp in A1
#  - updates country names + statistics for the rows whose country
#    order/data changed between the 10:50 and 11:20 snapshots

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 11:20"

# Row updates: country name (col A) + Casos totales/Nuevos casos/Casos activos/
# Recuperados/Casos criticos/Muertes hoy/Muertes (cols B-H)
$rows = @(
    @{ Row = 13; Name = 'Belgica'; Vals = @(12775,876,1527,10543,1021,192,705) },
    @{ Row = 17; Name = 'Austria'; Vals = @(9772,154,1095,8549,198,20,128) },
    @{ Row = 20; Name = 'Israel'; Vals = @(4831,136,163,4651,83,1,17) },
    @{ Row = 28; Name = 'Malasia'; Vals = @(2766,140,479,2244,94,6,43) },
    @{ Row = 55; Name = 'Eslovenia'; Vals = @(756,0,10,734,28,1,12) },
    @{ Row = 56; Name = 'Estonia'; Vals = @(745,30,26,715,13,1,4) },
    @{ Row = 57; Name = 'Hong Kong'; Vals = @(714,31,128,582,5,0,4) },
    @{ Row = 71; Name = 'Libano'; Vals = @(446,0,35,399,3,1,12) },
    @{ Row = 83; Name = 'Kuwait'; Vals = @(289,23,73,216,13,0,0) },
    @{ Row = 84; Name = 'Republica de Macedonia'; Vals = @(285,0,12,266,1,0,7) },
    @{ Row = 85; Name = 'Azerbaiyan'; Vals = @(273,0,26,243,23,0,4) },
    @{ Row = 86; Name = 'Jordania'; Vals = @(268,0,26,237,3,0,5) },
    @{ Row = 91; Name = 'Albania'; Vals = @(223,0,52,160,7,0,11) },
    @{ Row = 94; Name = 'Oman'; Vals = @(192,13,34,158,3,0,0) },
    @{ Row = 125; Name = 'Isla de Man'; Vals = @(60,11,0,60,0,0,0) },
    @{ Row = 126; Name = 'Aruba'; Vals = @(50,0,1,49,0,0,0) },
    @{ Row = 127; Name = 'Kenia'; Vals = @(50,0,1,48,2,0,1) },
    @{ Row = 146; Name = 'Etiopia'; Vals = @(25,2,4,21,1,0,0) },
    @{ Row = 183; Name = 'San Martin (Parte Holandesa)'; Vals = @(6,0,0,6,0,0,0) },
    @{ Row = 184; Name = 'Santa Sede'; Vals = @(6,0,0,6,0,0,0) },
    @{ Row = 185; Name = 'Benin'; Vals = @(6,0,1,5,0,0,0) },
    @{ Row = 190; Name = 'Islas Turcas y Caicos'; Vals = @(5,0,0,5,0,0,0) },
    @{ Row = 192; Name = 'Montserrat'; Vals = @(5,0,0,5,0,0,0) },
    @{ Row = 195; Name = 'Nicaragua'; Vals = @(4,0,0,3,0,0,1) },
    @{ Row = 196; Name = 'Gambia'; Vals = @(4,0,0,3,0,0,1) },
    @{ Row = 197; Name = 'Republica de Africa Central'; Vals = @(3,0,0,3,0,0,0) },
    @{ Row = 198; Name = 'Botsuana'; Vals = @(3,0,0,3,0,0,0) },
    @{ Row = 199; Name = 'Liberia'; Vals = @(3,0,0,3,0,0,0) },
    @{ Row = 200; Name = 'Islas Virgenes Britanicas'; Vals = @(3,1,0,3,0,0,0) },
    @{ Row = 201; Name = 'Belice'; Vals = @(3,0,0,3,0,0,0) }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.Name
    $col = 2
    foreach ($v in $r.Vals) {
        $ws.Cells.Item($rowNum, $col).Value = $v
        $col = $col + 1
    }
}
